$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (was LSPMW) becomes LSPM, keeping the 1/10, 25.6, 100 stats but with LSPM's own ECRPS value
$ws.Range("A5").Value = "LSPM"
$ws.Range("B5").Value = "1/10"
$ws.Range("C5").Value = 25.6
$ws.Range("D5").Value = 100
$ws.Range("E5").Value = 1.059559215525704

# Row 6 (was AV-MCPS) becomes LSPMW, taking on the 1/10, 25.6, 100 stats, but with LSPMW's own ECRPS value
$ws.Range("A6").Value = "LSPMW"
$ws.Range("B6").Value = "1/10"
$ws.Range("C6").Value = 25.6
$ws.Range("D6").Value = 100
$ws.Range("E6").Value = 1.056056548639601

# Row 9 (was LSPM) becomes AV-MCPS, taking on the 0/10, 0, 20 stats and AV-MCPS's own ECRPS value
$ws.Range("A9").Value = "AV-MCPS"
$ws.Range("B9").Value = "0/10"
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 20
$ws.Range("E9").Value = 2.740174974423968
